# Implement the "Medicine" card (and its companion LORD "The Healer")
# on the Cards sheet, rows 38-39, matching columns D:K.
#
# New shared strings must first be referenced in this order so they land
# at the expected indices: "The Healer", the Healer's ability text,
# Medicine's ability text, then "Medicine".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prime rows 38-39 with the same cell formatting used by the row above
# (row 37), so the new rows don't fall back to raw column defaults.
$ws.Range("E37:J37").Copy($ws.Range("E38:J38"))
$ws.Range("E37:J37").Copy($ws.Range("E39:J39"))

# Row 38 - Card ID 36: "The Healer" (LORD)
$ws.Range("E38").Value = "The Healer"
$ws.Range("F38").Value = "LORD"
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = "At the end of your turn, heal 1 unit you control to full health."
$ws.Range("K38").Value = "N"

# Row 39 - Card ID 37: "Medicine" (UTILITY)
$ws.Range("D39").Value = 37
$ws.Range("J39").Value = "Fully restore 1 unit's health."
$ws.Range("E39").Value = "Medicine"
$ws.Range("F39").Value = "UTILITY"
$ws.Range("G39").Value = 3
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = "Y"

# Restore the view: scrolled so column E is leftmost, with K39 selected.
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K39").Select()
